$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while preserving the numeric-looking
# literal as TEXT (not auto-converted to a Number by Excel), and restore the
# cells original Style so the quote-prefix flag added internally does not
# leak into the saved style.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "243.37"
Set-TextValue $ws.Range("D4") "5.404"
Set-TextValue $ws.Range("D6") "3.435"
Set-TextValue $ws.Range("D7") "6.531"
Set-TextValue $ws.Range("D8") "0.8094"
Set-TextValue $ws.Range("D9") "0.9236"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1425"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07424"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D12") "0.03309"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03087"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09355"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D15") "3.859"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001581"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D17") "0.04716"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0005932"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue $ws.Range("D19") "0.005868"
Set-TextValue $ws.Range("D20") "0.001277"
Set-TextValue $ws.Range("D21") "0.004891"
Set-TextValue $ws.Range("D22") "0.00006804"
Set-TextValue $ws.Range("D23") "3.569"
Set-TextValue $ws.Range("D40") "0.03971"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1078"
$ws.Range("E41").Value = "40BKEXTokenBKK"
Set-TextValue $ws.Range("D42") "0.002662"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.003070"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
Set-TextValue $ws.Range("D44") "0.009201"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
Set-TextValue $ws.Range("D45") "0.00005070"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue $ws.Range("D48") "0.002414"
Set-TextValue $ws.Range("D50") "0.0002000"
